$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

# Columns A-D hold text that Excel would otherwise auto-detect as date/number
# literals (e.g. "2024-01-11" -> date serial, "01" -> 1). Force text format,
# assign the literal value, then clear the formatting again so the cell ends
# up as a plain (unstyled) shared string, matching the rest of the sheet.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item($row, 1) "2024-01-11"
Set-TextValue $ws.Cells.Item($row, 2) "09:52:39"
$ws.Cells.Item($row, 3).Value = "Thursday"
Set-TextValue $ws.Cells.Item($row, 4) "01"

$ws.Cells.Item($row, 5).Value = 139423
$ws.Cells.Item($row, 6).Value = 142801
$ws.Cells.Item($row, 7).Value = 171208
$ws.Cells.Item($row, 8).Value = 148035
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119043
$ws.Cells.Item($row, 11).Value = 224688
$ws.Cells.Item($row, 12).Value = 251453
$ws.Cells.Item($row, 13).Value = 185360
$ws.Cells.Item($row, 14).Value = 110444
$ws.Cells.Item($row, 15).Value = 40769
$ws.Cells.Item($row, 16).Value = 30883
$ws.Cells.Item($row, 17).Value = 72795
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41814
$ws.Cells.Item($row, 20).Value = -1
